$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates: force text to avoid numeric auto-coercion ---
$priceCells = @{
    "D2" = "28.797.34"
    "D3" = "1.876.46"
    "D4" = "1.004"
    "D5" = "325.27"
    "D8" = "0.3870"
    "D9" = "0.07879"
    "D10" = "0.9872"
    "D12" = "1.904.13"
    "D13" = "7.006"
    "D14" = "5.707"
    "D15" = "0.06974"
    "D16" = "88.46"
    "D17" = "1.005"
    "D18" = "0.00001005"
    "D19" = "16.78"
    "D21" = "28.816.08"
    "D22" = "5.281"
    "D24" = "2.105"
    "D25" = "2.123.56"
    "D26" = "153.09"
    "D27" = "19.27"
    "D29" = "1.996"
    "D30" = "118.91"
    "D31" = "0.09339"
    "D32" = "0.9220"
    "D33" = "5.310"
    "D35" = "3.324"
    "D36" = "0.05797"
    "D37" = "1.149"
    "D38" = "0.02071"
    "D39" = "7.663"
    "D40" = "0.5635"
    "D42" = "9.790"
    "D43" = "0.07215"
    "D44" = "11.70"
    "D45" = "0.5303"
    "D46" = "2.151"
    "D48" = "1.841"
    "D49" = "113.40"
    "D51" = "1.004"
}
foreach ($ref in $priceCells.Keys) {
    $ws.Range($ref).NumberFormat = "@"
}
foreach ($ref in $priceCells.Keys) {
    $ws.Range($ref).Value = $priceCells[$ref]
}
foreach ($ref in $priceCells.Keys) {
    $ws.Range($ref).Style = "Normal"
}

# --- Column E (Volume 1h %) updates: plain text, safe as-is ---
$ws.Range("E2").Value = "  +2.51%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  +5.49%  "
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("E25").Value = "  +4.23%  "
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("E32").Value = "  -1.58%  "
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("E50").Value = "  +3.71%  "
$ws.Range("E51").Value = "  +0.32%  "
